# Generate Report for Handoff
#
# The localization status report previously recorded a failed handoff
# transform for 41b545a2-e135-46a8-afe1-f9a6297179d6.md. A handoff report
# has now been generated successfully for both the zh-cn and de-de
# targets, so:
#   * the "Handoff transform failed" status becomes "Ready for handoff"
#     (Overview summary + each language sheet)
#   * each language sheet gets a link to the newly produced handoff file
#   * the "Latest Handoff Datetime" is stamped with the generation time
#   * the "Handoff Reason" flips from "Ignored" to "Include" now that the
#     file is no longer being skipped

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTest/oltest/blob/b02f3e3848a83d70c4fb3eb84658976e950a945a"

# -- Overview sheet: roll-up status for both locales -----------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"

# -- Per-locale detail sheets ------------------------------------------------
$locales = @(
    @{ Sheet = "zh-cn"; HandoffFile = "41b545a2-e135-46a8-afe1-f9a6297179d6.33e66b2c014ddcaec4627ea9cdf090d4a36cab43.zh-cn.xlf"; HandoffTime = "2016-02-17 10:06:45" },
    @{ Sheet = "de-de"; HandoffFile = "41b545a2-e135-46a8-afe1-f9a6297179d6.33e66b2c014ddcaec4627ea9cdf090d4a36cab43.de-de.xlf"; HandoffTime = "2016-02-17 10:06:57" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    # Status: transform no longer fails, handoff is ready
    $ws.Range("B2").Value = "Ready for handoff"

    # Latest Handoff File: link to the freshly generated handoff package
    $ws.Hyperlinks.Add(
        $ws.Range("C2"),
        "$repoBase/$($locale.HandoffFile)",
        "",
        "",
        $locale.HandoffFile
    )

    # Latest Handoff Datetime: when the handoff file was generated
    $ws.Range("D2").Value = $locale.HandoffTime

    # Handoff Reason: file is now included in the handoff instead of ignored
    $ws.Range("H2").Value = "Include"
}
